$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC")

$ws.Range("D10").Value = 4

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D11").Select()

$wb.Save()
